{"js": "// Insert a new heading paragraph \"Initial post\" (bold, centered) at the\n// very start of the document body, ahead of all existing content.\nconst body = context.document.body;\n\nconst inserted = body.insertParagraph(\"Initial post\", Word.InsertLocation.start);\ninserted.alignment = Word.Alignment.centered;\ninserted.font.bold = true;\ninserted.font.boldBidirectional = true;\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$startRange = $d.Range(0, 0)\n$startRange.InsertParagraphBefore()\n\n$newPara = $d.Paragraphs(1).Range\n$newPara.Text = \"Initial post\"\n$newPara.ParagraphFormat.Alignment = 1  # wdAlignParagraphCenter\n$newPara.Bold = 1\n$newPara.BoldBi = 1\n"}
